$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 209.66
$ws.Range("I15").Value = 209.66
$ws.Range("K15").Value = 628.98
$ws.Range("M15").Value = -459.98
# Row 40
$ws.Range("H40").Value = 1049.3636
$ws.Range("I40").Value = 979.7
$ws.Range("J40").Value = 1107.4166
$ws.Range("K40").Value = 979.7
$ws.Range("L40").Value = 1107.4166
$ws.Range("M40").Value = -804.7
$ws.Range("N40").Value = -1457.4166
# Row 132
$ws.Range("H132").Value = 2748.4119
$ws.Range("I132").Value = 2749.0667
$ws.Range("K132").Value = 8247.2001
$ws.Range("M132").Value = -5717.2001
# Row 135
$ws.Range("H135").Value = 2757.9395
$ws.Range("I135").Value = 1675.8572
$ws.Range("J135").Value = 8817.6
$ws.Range("K135").Value = 15082.7148
$ws.Range("L135").Value = 79358.40000000001
$ws.Range("M135").Value = -12547.7148
$ws.Range("N135").Value = -84428.40000000001
# Row 139
$ws.Range("H139").Value = 70111.42999999999
$ws.Range("J139").Value = 70111.42999999999
$ws.Range("L139").Value = 70111.42999999999
$ws.Range("N139").Value = -80391.42999999999
# Row 140
$ws.Range("H140").Value = 87564
$ws.Range("J140").Value = 87564
$ws.Range("L140").Value = 87564
$ws.Range("N140").Value = -97924
# Row 141
$ws.Range("H141").Value = 2968.2812
$ws.Range("I141").Value = 3038.0435
$ws.Range("J141").Value = 2790
$ws.Range("K141").Value = 9114.130500000001
$ws.Range("L141").Value = 8370
$ws.Range("M141").Value = -3934.130500000001
$ws.Range("N141").Value = -18730

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 365489.5
$ws.Range("I32").Value = 5911.679
$ws.Range("K32").Value = 5911.679
$ws.Range("M32").Value = -5624.679
# Row 37
$ws.Range("H37").Value = 125007096
$ws.Range("J37").Value = 8112
$ws.Range("L37").Value = 8112
$ws.Range("N37").Value = -8658
# Row 57
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
# Row 61
$ws.Range("H61").Value = 3559.606
$ws.Range("I61").Value = 3595.7097
$ws.Range("K61").Value = 3595.7097
$ws.Range("M61").Value = -3383.7097
# Row 74
$ws.Range("H74").Value = 1159.4814
$ws.Range("I74").Value = 694.6
$ws.Range("J74").Value = 2487.7144
$ws.Range("K74").Value = 694.6
$ws.Range("L74").Value = 2487.7144
$ws.Range("M74").Value = 179.4
$ws.Range("N74").Value = -4235.7144
# Row 77
$ws.Range("H77").Value = 1159.4814
$ws.Range("I77").Value = 694.6
$ws.Range("J77").Value = 2487.7144
$ws.Range("K77").Value = 3473
$ws.Range("L77").Value = 12438.572
$ws.Range("M77").Value = 895
$ws.Range("N77").Value = -21174.572
# Row 132
$ws.Range("H132").Value = 2077.6
$ws.Range("I132").Value = 1770.5416
$ws.Range("K132").Value = 5311.6248
$ws.Range("M132").Value = -2781.6248
# Row 136
$ws.Range("H136").Value = 3559.606
$ws.Range("I136").Value = 3595.7097
$ws.Range("K136").Value = 10787.1291
$ws.Range("M136").Value = -8237.1291

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1415.6428
$ws.Range("I134").Value = 878
$ws.Range("J134").Value = 2759.75
$ws.Range("K134").Value = 2634
$ws.Range("L134").Value = 8279.25
$ws.Range("M134").Value = -99
$ws.Range("N134").Value = -13349.25
# Row 140
$ws.Range("H140").Value = 89700
$ws.Range("J140").Value = 89700
$ws.Range("L140").Value = 89700
$ws.Range("N140").Value = -100060

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 50
$ws.Range("H50").Value = 8914
$ws.Range("J50").Value = 8914
$ws.Range("L50").Value = 8914
$ws.Range("N50").Value = -10164
# Row 51
$ws.Range("H51").Value = 8380.6
$ws.Range("I51").Value = 4500
$ws.Range("J51").Value = 9350.75
$ws.Range("K51").Value = 4500
$ws.Range("L51").Value = 9350.75
$ws.Range("M51").Value = -3764
$ws.Range("N51").Value = -10822.75
# Row 59
$ws.Range("H59").Value = 11927
$ws.Range("J59").Value = 11927
$ws.Range("L59").Value = 11927
$ws.Range("N59").Value = -14217
# Row 60
$ws.Range("H60").Value = 8251
$ws.Range("J60").Value = 8251
$ws.Range("L60").Value = 8251
$ws.Range("N60").Value = -9273
# Row 61
$ws.Range("H61").Value = 8380.6
$ws.Range("I61").Value = 4500
$ws.Range("J61").Value = 9350.75
$ws.Range("K61").Value = 4500
$ws.Range("L61").Value = 9350.75
$ws.Range("M61").Value = -4152
$ws.Range("N61").Value = -10046.75
# Row 74
$ws.Range("H74").Value = 15037.429
$ws.Range("J74").Value = 17162.834
$ws.Range("L74").Value = 17162.834
$ws.Range("N74").Value = -18910.834
# Row 77
$ws.Range("H77").Value = 15037.429
$ws.Range("J77").Value = 17162.834
$ws.Range("L77").Value = 51488.50199999999
$ws.Range("N77").Value = -60224.50199999999
# Row 132
$ws.Range("H132").Value = 3430.1428
$ws.Range("I132").Value = 2303.4285
$ws.Range("J132").Value = 4556.857
$ws.Range("K132").Value = 6910.2855
$ws.Range("L132").Value = 13670.571
$ws.Range("M132").Value = -4380.2855
$ws.Range("N132").Value = -18730.571
# Row 140
$ws.Range("H140").Value = 90000
$ws.Range("J140").Value = 90000
$ws.Range("L140").Value = 90000
$ws.Range("N140").Value = -100360

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Range("H38").Value = 472.8125
$ws.Range("I38").Value = 888
$ws.Range("J38").Value = 223.7
$ws.Range("K38").Value = 2664
$ws.Range("L38").Value = 671.0999999999999
$ws.Range("M38").Value = -2317
$ws.Range("N38").Value = -1365.1
# Row 81
$ws.Range("H81").Value = 2526.625
$ws.Range("I81").Value = 913
$ws.Range("J81").Value = 2757.1428
$ws.Range("K81").Value = 2739
$ws.Range("L81").Value = 8271.428400000001
$ws.Range("M81").Value = -1616
$ws.Range("N81").Value = -10517.4284
# Row 84
$ws.Range("H84").Value = 2526.625
$ws.Range("I84").Value = 913
$ws.Range("J84").Value = 2757.1428
$ws.Range("K84").Value = 8217
$ws.Range("L84").Value = 24814.2852
$ws.Range("M84").Value = -2601
$ws.Range("N84").Value = -36046.2852

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 82
$ws.Range("H82").Value = 33960
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 33960
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 33960
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -34726
# Row 85
$ws.Range("H85").Value = 33960
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 33960
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 33960
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -36612
# Row 132
$ws.Range("H132").Value = 2580.862
$ws.Range("I132").Value = 2034.15
$ws.Range("J132").Value = 3795.7778
$ws.Range("K132").Value = 6102.450000000001
$ws.Range("L132").Value = 11387.3334
$ws.Range("M132").Value = -3572.450000000001
$ws.Range("N132").Value = -16447.3334
# Row 140
$ws.Range("H140").Value = 99879.664
$ws.Range("J140").Value = 99879.664
$ws.Range("L140").Value = 99879.664
$ws.Range("N140").Value = -110239.664

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 3289
$ws.Range("I61").Value = 1946.8334
$ws.Range("J61").Value = 4631.1665
$ws.Range("K61").Value = 1946.8334
$ws.Range("L61").Value = 4631.1665
$ws.Range("M61").Value = -1744.8334
$ws.Range("N61").Value = -5035.1665
# Row 113
$ws.Range("H113").Value = 3289
$ws.Range("I113").Value = 1946.8334
$ws.Range("J113").Value = 4631.1665
$ws.Range("K113").Value = 1946.8334
$ws.Range("L113").Value = 4631.1665
$ws.Range("M113").Value = 223.1666
$ws.Range("N113").Value = -8971.166499999999
# Row 122
$ws.Range("H122").Value = 2776.75
$ws.Range("I122").Value = 2313.4546
$ws.Range("J122").Value = 3796
$ws.Range("K122").Value = 6940.3638
$ws.Range("L122").Value = 11388
$ws.Range("M122").Value = -4490.3638
$ws.Range("N122").Value = -16288
# Row 132
$ws.Range("H132").Value = 3300.0908
$ws.Range("I132").Value = 3443.44
$ws.Range("J132").Value = 2852.125
$ws.Range("K132").Value = 10330.32
$ws.Range("L132").Value = 8556.375
$ws.Range("M132").Value = -7800.32
$ws.Range("N132").Value = -13616.375
# Row 136
$ws.Range("H136").Value = 4290.8213
$ws.Range("I136").Value = 3000.5334
$ws.Range("J136").Value = 5779.615
$ws.Range("K136").Value = 9001.600199999999
$ws.Range("L136").Value = 17338.845
$ws.Range("M136").Value = -6451.600199999999
$ws.Range("N136").Value = -22438.845

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 725.5484
$ws.Range("I113").Value = 538.36
$ws.Range("J113").Value = 1505.5
$ws.Range("K113").Value = 1615.08
$ws.Range("L113").Value = 4516.5
$ws.Range("M113").Value = 554.9200000000001
$ws.Range("N113").Value = -8856.5
# Row 122
$ws.Range("H122").Value = 1765.1613
$ws.Range("I122").Value = 1513.0435
$ws.Range("J122").Value = 2490
$ws.Range("K122").Value = 4539.1305
$ws.Range("L122").Value = 7470
$ws.Range("M122").Value = -2089.1305
$ws.Range("N122").Value = -12370
# Row 132
$ws.Range("H132").Value = 34092476
$ws.Range("I132").Value = 51725450
$ws.Range("K132").Value = 155176350
$ws.Range("M132").Value = -155173820
# Row 136
$ws.Range("H136").Value = 640.35596
$ws.Range("I136").Value = 615.39215
$ws.Range("K136").Value = 1846.17645
$ws.Range("M136").Value = 703.8235500000001
